$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header: student name (text unchanged, but keep consistent) ---
$ws.Cells.Item(4, 2).Value = "Amitrajit Sarkar"

# --- Row 7: refine the list item style ---
$ws.Cells.Item(7, 1).Value = 43575
$ws.Cells.Item(7, 7).Value = "refine the list item sytle"
$ws.Cells.Item(7, 8).Value = "The guess history display too large"
$ws.Rows.Item(7).RowHeight = 16

# --- Row 8: refine the Start button style ---
$ws.Cells.Item(8, 1).Value = 43575
$ws.Cells.Item(8, 7).Value = "refine the Start button style"
$ws.Cells.Item(8, 8).Value = "The Start button is too small"

# --- Row 9: adjust the feedback button align ---
$ws.Cells.Item(9, 1).Value = 43575
$ws.Cells.Item(9, 7).Value = "adjust the feedback button align"
$ws.Cells.Item(9, 8).Value = "The feedback button isn't align good"

# --- Row 10: extend the padding ---
$ws.Cells.Item(10, 1).Value = 43575
$ws.Cells.Item(10, 7).Value = "extend the padding"
$ws.Cells.Item(10, 8).Value = "The padding under the header  is too narrow"

# --- Row 11: refine the margin ---
$ws.Cells.Item(11, 1).Value = 43575
$ws.Cells.Item(11, 7).Value = "refine the margin"
$ws.Cells.Item(11, 8).Value = "The margin of the application is not good"

# --- View state: active selection moves to A14 ---
$ws.Activate() | Out-Null
$ws.Range("A14").Select() | Out-Null
